$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.42%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.62"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.89%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.009"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.45%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07855"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.21%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.220"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.90%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.023"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.03%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.005"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.82%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9094"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.12%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1862"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.67%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09228"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-9.21%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08466"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.54%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03524"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.30%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09947"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.52%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001469"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.87%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005657"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.90%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.469"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.05%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.153"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.72%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.87%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.18%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.794"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "11.03%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2201"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-7.77%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04646"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.77%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001227"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.91%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.35%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001297"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.22%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004744"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "39.82%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01761"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.05%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04732"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.50%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007869"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.47%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1393"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.36%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007657"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.61%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002215"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.02%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01023"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "11.22%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.10%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.10%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.670"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "216.29%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "34.35%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.10%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.10%"

Write-Output "Applied cryptos.xlsx price/volume updates"
